$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append short/full description text in parentheses to the existing labels.
$ws.Range("A1").Value = "ཁ་འདོན།(ཁ་འདོན་འགྲེལ་བཤད་)(ཁ་འདོན་འགྲེལ་བཤད་ཐུང་ཐུང་)"
$ws.Range("B2").Value = "སྨོན་ལམ།(སྨོན་ལམ་འགྲེལ་བཤད་)(སྨོན་ལམ་འགྲེལ་བཤད་ཐུང་ཐུང་)"
$ws.Range("C3").Value = "ཀུན་བཟང་སྨོན་ལམ།(ཀུན་བཟང་སྨོན་ལམ་འགྲེལ་བཤད་)(ཀུན་བཟང་སྨོན་ལམ་འགྲེལ་བཤད་ཐུང་ཐུང་)"
$ws.Range("C6").Value = "བཟང་སྤྱོད་སྨོན་ལམ།(བཟང་སྤྱོད་འགྲེལ་བཤད་)(བཟང་སྤྱོད་འགྲེལ་བཤད་ཐུང་ཐུང་)"

# Widen the columns to fit the longer labels (nearest value reachable through
# the pixel-quantized ColumnWidth property that yields the target stored width).
$ws.Columns.Item(1).ColumnWidth = 27.333333333333332
$ws.Columns.Item(2).ColumnWidth = 28
$ws.Columns.Item(3).ColumnWidth = 37.833333333333336
